$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 2 ("Shana AFFICHARD"),
# pushing it (and the rows below) down.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Insert two new rows before the current row 5 ("Anna ANDREY"),
# pushing it down.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Row 2: Ina AÏSSI
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Ina AÏSSI"
$ws.Range("C2").Value = "5 RIO GRANDE"
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = "F"
$ws.Range("F2").Value = "00:11:34"

# Row 3: Robin BEAUFILS
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Robin BEAUFILS"
$ws.Range("C3").Value = "5 AMAZONE"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "G"
$ws.Range("F3").Value = "00:11:34"

# Row 4: Shana AFFICHARD (existing row, update rank only)
$ws.Range("A4").Value = 3

# Row 5: Yanis ALLIX
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Yanis ALLIX"
$ws.Range("C5").Value = "5 RIO GRANDE"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = "G"
$ws.Range("F5").Value = "00:11:34"

# Row 6: Thaïs BAILLARD
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Thaïs BAILLARD"
$ws.Range("C6").Value = "5 SANTA CRUZ"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = "F"
$ws.Range("F6").Value = "00:11:34"

# Row 7: Anna ANDREY (existing row, update rank only)
$ws.Range("A7").Value = 6
